# Apply the 8.6.1 indicator update:
#  - fix the Kyrgyz NEET indicator title in row 1 (A1) and its row height
#  - add the 2023 data column (T) mirroring the 2022 column (S) formatting
#  - move the sheet selection back to A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: corrected Kyrgyz title text, shorter row height ---
$ws.Range("A1").Value = "8.6.1 Иштебеген, окубаган жана кесиптик көндүмдөрдү үйрөнбөгөн (15 жаштан 24 жашка чейинки ) жаштардын үлүшү  "
$ws.Range("A1").VerticalAlignment = -4108   # xlCenter
$ws.Rows.Item(1).RowHeight = 48

# --- New column T: 2023 data, copying formatting from column S ---
$ws.Range("S4:S7").Copy($ws.Range("T4:T7"))
$ws.Range("T4").Value = 2023
$ws.Range("T5").Value = 18.6
$ws.Range("T6").Value = 11.5
$ws.Range("T7").Value = 25.9

# --- Reset the active selection to A1 (default view) ---
$ws.Range("A1").Select() | Out-Null
